$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.318.19"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.671.02"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.07"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.49"
$ws.Range("E6").Value = "  -6.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.668.99"
$ws.Range("E7").Value = "  -3.74%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -5.23%  "
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.41"
$ws.Range("E13").Value = "  -6.16%  "
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.281.31"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.668.86"
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.377.69"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.66"
$ws.Range("E18").Value = "  +6.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.14"
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.47"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.10"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.32"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -6.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.12"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.35"
$ws.Range("E31").Value = "  -6.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.34"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.807.77"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  -5.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.607.34"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.987"
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.74"
$ws.Range("E39").Value = "  -5.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.321"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.13"
$ws.Range("E42").Value = "  -10.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.58"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -5.97%  "
$ws.Range("E45").Value = "  -8.40%  "
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.42"
$ws.Range("E48").Value = "  -7.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.67"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.747.62"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("E51").Value = "  -3.88%  "
